# Apply the edit described by the diff:
# 1. Rename worksheet "HP_Test" -> "HP"
# 2. Append a new data row (row 6) with serial_number/brand/modele/product_number,
#    reusing the existing formatting of column A (bold + border + centered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "HP"

# Copy the format of A5 (bold/border/centered) down to the new A6 cell first,
# so the new row reuses the existing style instead of creating a new one.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats

# Add the new row of data
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "ACM029T0XJ"
$ws.Range("C6").Value = "HPE"
$ws.Range("D6").Value = "MSA 1050"
$ws.Range("E6").Value = "Q2R21B"
